# Adding data driven approach from excel
# Extends the TestData sheet with new columns (Language, PhoneNo, Password,
# Customer Name, Customer Ph, Credit, Payment) and a second data row (TC_02).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row (row 1): new column headers E1:H1, styled like A1:D1
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Language"
$ws.Range("C1").Value = "PhoneNo"
$ws.Range("D1").Value = "Password"

$ws.Range("A1").Copy()
$ws.Range("E1:H1").PasteSpecial(-4122)

$ws.Range("E1").Value = "Customer Name"
$ws.Range("F1").Value = "Customer Ph"
$ws.Range("G1").Value = "Credit"
$ws.Range("H1").Value = "Payment"

# ---------------------------------------------------------------------
# 2. Row 2 (existing TC_01 row): fill in the new columns B:D and E,F,G,H
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "TC_01"
$ws.Range("B2").Value = "English"
$ws.Range("C2").Value = "7406764431"
$ws.Range("D2").Value = "123456"

$e2 = $ws.Range("E2")
$e2.Value = "Test User 1"
$e2.Font.Name = "Calibri"
$e2.NumberFormat = "@"
$e2.HorizontalAlignment = -4131

$f2 = $ws.Range("F2")
$f2.Font.Name = "Calibri"
$f2.NumberFormat = "@"

$g2 = $ws.Range("G2")
$g2.Value = "250"
$g2.Font.Name = "Calibri"
$g2.NumberFormat = "@"
$g2.HorizontalAlignment = -4131

$h2 = $ws.Range("H2")
$h2.Value = "150"
$h2.Font.Name = "Calibri"
$h2.NumberFormat = "@"
$h2.HorizontalAlignment = -4131

# ---------------------------------------------------------------------
# 3. Row 3 (new TC_02 row): same shape as row 2
# ---------------------------------------------------------------------
$ws.Range("A2:D2").Copy()
$ws.Range("A3:D3").PasteSpecial(-4122)

$ws.Range("A3").Value = "TC_02"
$ws.Range("B3").Value = "English"
$ws.Range("C3").Value = "7406764431"
$ws.Range("D3").Value = "123456"

$e3 = $ws.Range("E3")
$e3.Value = "Test User 2"
$e3.Font.Name = "Calibri"
$e3.NumberFormat = "@"
$e3.HorizontalAlignment = -4131

$f3 = $ws.Range("F3")
$f3.Font.Name = "Calibri"
$f3.NumberFormat = "@"

$g3 = $ws.Range("G3")
$g3.Value = "250"
$g3.Font.Name = "Calibri"
$g3.NumberFormat = "@"
$g3.HorizontalAlignment = -4131

$h3 = $ws.Range("H3")
$h3.Value = "250"
$h3.Font.Name = "Calibri"
$h3.NumberFormat = "@"
$h3.HorizontalAlignment = -4131

# Fill the remainder of row 3 (I3:Z3) with the same blank style used on
# the rest of the sheet (I1:Z1 / I2:Z2).
$ws.Range("I1:Z1").Copy()
$ws.Range("I3:Z3").PasteSpecial(-4122)
